$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 89 mirrors the formatting of row 88 (header-free data row):
# column A bold/bordered/centered style, column E date-time number format,
# all the other columns left with the default (unstyled) format.
$ws.Range("A88:V88").Copy()
$ws.Range("A89:V89").PasteSpecial(-4122)

# Plain text / label columns
$ws.Range("B89").Value = "ecuador"
$ws.Range("C89").Value = "liga-pro"

# D89 must stay a text string ("2023"), not auto-convert to a number;
# the leading apostrophe forces text, then reapplying the Normal style
# clears the quote-prefix formatting iron_native adds for that.
$ws.Range("D89").Value = "'2023"
$ws.Range("D89").Style = "Normal"

# Numeric index / score / odds columns
$ws.Range("A89").Value = 88
$ws.Range("E89").Value = 45234.04166666666
$ws.Range("F89").Value = "Dep. Cuenca"
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = "Aucas"
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2.16
$ws.Range("K89").Value = "29/10/2023 21:42"
$ws.Range("L89").Value = 2.74
$ws.Range("M89").Value = "04/11/2023 00:52"
$ws.Range("N89").Value = 3.37
$ws.Range("O89").Value = "29/10/2023 21:42"
$ws.Range("P89").Value = 3.34
$ws.Range("Q89").Value = "04/11/2023 00:50"
$ws.Range("R89").Value = 3.4
$ws.Range("S89").Value = "29/10/2023 21:42"
$ws.Range("T89").Value = 2.64
$ws.Range("U89").Value = "04/11/2023 00:52"
$ws.Range("V89").Value = "https://www.betexplorer.com/football/ecuador/liga-pro/dep-cuenca-aucas/8lMTJSkB/"
